# Updated cryptos list on Wed May 31 08:25:49 UTC 2023 with GitHub Actions
#
# Refresh the Coin / Link / Price / Volume(1h) table on the active sheet
# with the latest scraped coinranking.com values. Rows 12/13 and 37/38
# also swap their Coin+Link pairs (ranking order changed between runs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be read back as plain text rather than being
# auto-coerced to a number by the COM value setter (values like "1.000",
# "27.199.36" or "1.170" would otherwise lose their literal formatting).
$ws.Range("D2:D51").NumberFormat = "@"

# --- Volume(1h) (column E) updates ---
$ws.Range("E3").Value  = "  -1.80%  "
$ws.Range("E4").Value  = "  -0.15%  "
$ws.Range("E5").Value  = "  -1.84%  "
$ws.Range("E6").Value  = "  -0.22%  "
$ws.Range("E7").Value  = "  +1.03%  "
$ws.Range("E8").Value  = "  -1.24%  "
$ws.Range("E9").Value  = "  -1.53%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  -0.42%  "

# --- Row 12 / Row 13 swap to WrappedEther / TRON ---
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E13").Value = "  -1.07%  "

$ws.Range("E14").Value = "  -2.69%  "
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("E22").Value = "  -1.55%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("E30").Value = "  -2.79%  "
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  -1.57%  "
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("E36").Value = "  -4.04%  "

# --- Row 37 / Row 38 swap to RenderToken / VeChain ---
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E37").Value = "  +2.39%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E38").Value = "  -1.61%  "

$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("E41").Value = "  -2.84%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("E51").Value = "  -0.35%  "

# --- Price (column D) updates ---
$ws.Range("D2").Value  = "27.199.36"
$ws.Range("D3").Value  = "1.872.40"
$ws.Range("D5").Value  = "307.63"
$ws.Range("D6").Value  = "1.000"
$ws.Range("D7").Value  = "0.5082"
$ws.Range("D9").Value  = "0.07162"
$ws.Range("D10").Value = "0.8921"
$ws.Range("D11").Value = "20.76"
$ws.Range("D12").Value = "1.894.14"
$ws.Range("D13").Value = "0.07578"
$ws.Range("D14").Value = "5.337"
$ws.Range("D15").Value = "89.46"
$ws.Range("D17").Value = "0.000008547"
$ws.Range("D18").Value = "14.16"
$ws.Range("D19").Value = "1.000"
$ws.Range("D20").Value = "27.238.38"
$ws.Range("D21").Value = "5.078"
$ws.Range("D22").Value = "2.107.25"
$ws.Range("D24").Value = "6.507"
$ws.Range("D25").Value = "150.64"
$ws.Range("D26").Value = "1.846"
$ws.Range("D27").Value = "18.01"
$ws.Range("D28").Value = "2.118"
$ws.Range("D30").Value = "4.763"
$ws.Range("D31").Value = "4.728"
$ws.Range("D32").Value = "0.08995"
$ws.Range("D34").Value = "3.097"
$ws.Range("D35").Value = "0.7533"
$ws.Range("D36").Value = "1.170"
$ws.Range("D37").Value = "2.563"
$ws.Range("D38").Value = "0.02031"
$ws.Range("D39").Value = "3.043"
$ws.Range("D40").Value = "1.076"
$ws.Range("D41").Value = "0.5366"
$ws.Range("D42").Value = "6.626"
$ws.Range("D43").Value = "114.68"
$ws.Range("D44").Value = "8.516"
$ws.Range("D45").Value = "0.1485"
$ws.Range("D46").Value = "0.4673"
$ws.Range("D47").Value = "1.000"
$ws.Range("D48").Value = "10.11"
$ws.Range("D49").Value = "1.572"
$ws.Range("D50").Value = "65.16"
$ws.Range("D51").Value = "36.77"

# Drop the temporary text-format stamp back to the sheet's default style
# so the saved cells carry no extra formatting vs. before.xlsx.
$ws.Range("D2:D51").Style = "Normal"

Write-Output "cryptos refreshed"
